$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.328.18"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "1.804.12"

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.574"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.72%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.89%  "

# Row 9
$ws.Range("E9").Value = "  +2.10%  "

# Row 10
$ws.Range("E10").Value = "  +0.45%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0966"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.06%  "

# Row 12
$ws.Range("D12").Value = "2.065.31"
$ws.Range("E12").Value = "  +0.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.36%  "

# Row 14
$ws.Range("D14").Value = "1.804.13"
$ws.Range("E14").Value = "  +0.68%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.40%  "

# Row 17
$ws.Range("D17").Value = "34.312.23"
$ws.Range("E17").Value = "  -0.11%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20
$ws.Range("E20").Value = "  -0.17%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.94%  "

# Row 22
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.67%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "172.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.73%  "

# Row 25
$ws.Range("E25").Value = "  +3.22%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.51%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.97%  "

# Row 28
$ws.Range("E28").Value = "  +2.67%  "

# Row 29
$ws.Range("E29").Value = "  -0.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0532"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.32%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.41%  "

# Row 35
$ws.Range("D35").Value = "1.392.87"
$ws.Range("E35").Value = "  -1.23%  "

# Row 36
$ws.Range("E36").Value = "  -1.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.22%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.61%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "

# Row 40
$ws.Range("E40").Value = "  +10.62%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.963"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.53%  "

# Row 42
$ws.Range("E42").Value = "  +1.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "81.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.69%  "

# Row 44
$ws.Range("E44").Value = "  +0.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.88%  "

# Row 46
$ws.Range("E46").Value = "  -0.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0502"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.83%  "

# Row 48
$ws.Range("D48").Value = "1.965.34"
$ws.Range("E48").Value = "  +0.92%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.73"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "  +0.05%  "

# Row 51
$ws.Range("E51").Value = "  +0.13%  "
